$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2252.296
$ws.Range("I15").Value = 2252.296
$ws.Range("K15").Value = 6756.887999999999
$ws.Range("M15").Value = -6587.887999999999

$ws.Range("H17").Value = 177.76086
$ws.Range("J17").Value = 172.72093
$ws.Range("L17").Value = 518.1627900000001
$ws.Range("N17").Value = -854.1627900000001

$ws.Range("H112").Value = 6631.5386
$ws.Range("J112").Value = 7055
$ws.Range("L112").Value = 21165
$ws.Range("N112").Value = -23381

$ws.Range("H137").Value = 1211.4595
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1211.4595
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 3634.3785
$ws.Range("N137").Value = -8734.378499999999
$ws.Range("M137").Value = ""

$ws.Range("H138").Value = 2756.4167
$ws.Range("I138").Value = 5220.5713
$ws.Range("J138").Value = 2430.9622
$ws.Range("K138").Value = 15661.7139
$ws.Range("L138").Value = 7292.8866
$ws.Range("M138").Value = -10521.7139
$ws.Range("N138").Value = -17572.8866

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 727910.4
$ws.Range("I32").Value = 895343.4399999999
$ws.Range("J32").Value = 16320.125
$ws.Range("K32").Value = 895343.4399999999
$ws.Range("L32").Value = 16320.125
$ws.Range("M32").Value = -895056.4399999999
$ws.Range("N32").Value = -16894.125

$ws.Range("H61").Value = 2855.889
$ws.Range("I61").Value = 1200.6
$ws.Range("J61").Value = 4925
$ws.Range("K61").Value = 1200.6
$ws.Range("L61").Value = 4925
$ws.Range("M61").Value = -988.5999999999999
$ws.Range("N61").Value = -5349

$ws.Range("H74").Value = 1825.4375
$ws.Range("I74").Value = 1273.5714
$ws.Range("J74").Value = 2254.6667
$ws.Range("K74").Value = 1273.5714
$ws.Range("L74").Value = 2254.6667
$ws.Range("M74").Value = -399.5714
$ws.Range("N74").Value = -4002.6667

$ws.Range("H77").Value = 1825.4375
$ws.Range("I77").Value = 1273.5714
$ws.Range("J77").Value = 2254.6667
$ws.Range("K77").Value = 6367.857
$ws.Range("L77").Value = 11273.3335
$ws.Range("M77").Value = -1999.857
$ws.Range("N77").Value = -20009.3335

$ws.Range("H132").Value = 4829.231
$ws.Range("I132").Value = 4356.6206
$ws.Range("J132").Value = 6199.8
$ws.Range("K132").Value = 13069.8618
$ws.Range("L132").Value = 18599.4
$ws.Range("M132").Value = -10539.8618
$ws.Range("N132").Value = -23659.4

$ws.Range("H136").Value = 2855.889
$ws.Range("I136").Value = 1200.6
$ws.Range("J136").Value = 4925
$ws.Range("K136").Value = 3601.8
$ws.Range("L136").Value = 14775
$ws.Range("M136").Value = -1051.8
$ws.Range("N136").Value = -19875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 79065
$ws.Range("J45").Value = 79065
$ws.Range("L45").Value = 79065
$ws.Range("N45").Value = -80681

$ws.Range("H105").Value = 6947039.5
$ws.Range("I105").Value = 9617516
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 9617516
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -9615769
$ws.Range("N105").Value = -7294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3503.196
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3503.196
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3503.196
$ws.Range("N31").Value = -4093.196
$ws.Range("M31").Value = ""

$ws.Range("H34").Value = 3503.196
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3503.196
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3503.196
$ws.Range("N34").Value = -3907.196
$ws.Range("M34").Value = ""

$ws.Range("H107").Value = 3472898.8
$ws.Range("I107").Value = 5682352.5
$ws.Range("K107").Value = 5682352.5
$ws.Range("M107").Value = -5680432.5

$ws.Range("H122").Value = 1701.3462
$ws.Range("I122").Value = 1223.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3670.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1220.5
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1185.881
$ws.Range("I5").Value = 537.3125
$ws.Range("J5").Value = 1585
$ws.Range("K5").Value = 1611.9375
$ws.Range("L5").Value = 4755
$ws.Range("M5").Value = -1499.9375
$ws.Range("N5").Value = -4979

$ws.Range("H68").Value = 3291.446
$ws.Range("I68").Value = 7166.467
$ws.Range("J68").Value = 2128.94
$ws.Range("K68").Value = 21499.401
$ws.Range("L68").Value = 6386.82
$ws.Range("M68").Value = -20688.401
$ws.Range("N68").Value = -8008.82

$ws.Range("H71").Value = 3291.446
$ws.Range("I71").Value = 7166.467
$ws.Range("J71").Value = 2128.94
$ws.Range("K71").Value = 64498.20299999999
$ws.Range("L71").Value = 19160.46
$ws.Range("M71").Value = -60442.20299999999
$ws.Range("N71").Value = -27272.46

$ws.Range("H74").Value = 4000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -14122
$ws.Range("M74").Value = ""

$ws.Range("H77").Value = 4000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -46608
$ws.Range("M77").Value = ""

$ws.Range("H81").Value = 7167.643
$ws.Range("I81").Value = 1004.3333
$ws.Range("J81").Value = 8848.546
$ws.Range("K81").Value = 3012.9999
$ws.Range("L81").Value = 26545.638
$ws.Range("M81").Value = -1889.9999
$ws.Range("N81").Value = -28791.638

$ws.Range("H84").Value = 7167.643
$ws.Range("I84").Value = 1004.3333
$ws.Range("J84").Value = 8848.546
$ws.Range("K84").Value = 9038.9997
$ws.Range("L84").Value = 79636.914
$ws.Range("M84").Value = -3422.9997
$ws.Range("N84").Value = -90868.914

$ws.Range("H107").Value = 2213.8616
$ws.Range("I107").Value = 326.90475
$ws.Range("J107").Value = 3114.4546
$ws.Range("K107").Value = 980.71425
$ws.Range("L107").Value = 9343.363799999999
$ws.Range("M107").Value = 939.28575
$ws.Range("N107").Value = -13183.3638

$ws.Range("H131").Value = 730.5
$ws.Range("I131").Value = 321.46155
$ws.Range("J131").Value = 1085
$ws.Range("K131").Value = 964.38465
$ws.Range("L131").Value = 3255
$ws.Range("M131").Value = 4075.61535
$ws.Range("N131").Value = -13335

$ws.Range("H132").Value = 2858.5432
$ws.Range("J132").Value = 3801.0732
$ws.Range("L132").Value = 34209.6588
$ws.Range("N132").Value = -39269.6588

$ws.Range("H135").Value = 1185.881
$ws.Range("I135").Value = 537.3125
$ws.Range("J135").Value = 1585
$ws.Range("K135").Value = 4835.8125
$ws.Range("L135").Value = 14265
$ws.Range("M135").Value = -2300.8125
$ws.Range("N135").Value = -19335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = ""

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""

$ws.Range("H113").Value = 1512.6666
$ws.Range("I113").Value = 922.2
$ws.Range("K113").Value = 922.2
$ws.Range("M113").Value = 1247.8

$ws.Range("H122").Value = 5212.8184
$ws.Range("I122").Value = 4165.6665
$ws.Range("J122").Value = 5605.5
$ws.Range("K122").Value = 12496.9995
$ws.Range("L122").Value = 16816.5
$ws.Range("M122").Value = -10046.9995
$ws.Range("N122").Value = -21716.5

$ws.Range("H141").Value = 62464.5
$ws.Range("J141").Value = 62464.5
$ws.Range("L141").Value = 62464.5
$ws.Range("N141").Value = -72824.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2760.1887
$ws.Range("I132").Value = 2481.946
$ws.Range("J132").Value = 3403.625
$ws.Range("K132").Value = 7445.838
$ws.Range("L132").Value = 10210.875
$ws.Range("M132").Value = -4915.838
$ws.Range("N132").Value = -15270.875

$ws.Range("H136").Value = 3031284
$ws.Range("I136").Value = 994.0968
$ws.Range("K136").Value = 2982.2904
$ws.Range("M136").Value = -432.2903999999999

$ws.Range("H137").Value = 40276.332
$ws.Range("J137").Value = 40276.332
$ws.Range("L137").Value = 40276.332
$ws.Range("N137").Value = -50476.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 23333.334
$ws.Range("J51").Value = 23333.334
$ws.Range("L51").Value = 23333.334
$ws.Range("N51").Value = -24353.334

$ws.Range("H122").Value = 2664.3635
$ws.Range("I122").Value = 2272.5715
$ws.Range("J122").Value = 3350
$ws.Range("K122").Value = 6817.7145
$ws.Range("L122").Value = 10050
$ws.Range("M122").Value = -4367.7145
$ws.Range("N122").Value = -14950

$ws.Range("H123").Value = 35627.375
$ws.Range("J123").Value = 35627.375
$ws.Range("L123").Value = 35627.375
$ws.Range("N123").Value = -45427.375

$ws.Range("H141").Value = 69566.664
$ws.Range("J141").Value = 69566.664
$ws.Range("L141").Value = 69566.664
$ws.Range("N141").Value = -79926.664
